$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells that look like plain numbers stay as text,
# matching the original inlineStr string cells (preserve "582.20" style formatting).
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D11", "D13", "D14", "D16", "D17", "D19", "D20", "D21", "D22", "D24", "D25", "D28", "D30", "D34", "D35", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D47", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "63.428.36"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "3.092.65"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "582.20"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").Value = "144.67"
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.086.16"
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("E10").Value = "  +6.16%  "
$ws.Range("D11").Value = "5.59"
$ws.Range("E11").Value = "  -2.87%  "
$ws.Range("E12").Value = "  -2.63%  "
$ws.Range("D13").Value = "0.0000245"
$ws.Range("E13").Value = "  -1.69%  "
$ws.Range("D14").Value = "37.22"
$ws.Range("E14").Value = "  +4.16%  "
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").Value = "3.611.87"
$ws.Range("D17").Value = "63.339.08"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("E18").Value = "  -1.34%  "
$ws.Range("D19").Value = "3.094.72"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").Value = "461.58"
$ws.Range("E20").Value = "  -1.54%  "
$ws.Range("D21").Value = "14.22"
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("D22").Value = "0.723"
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("D24").Value = "81.23"
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("D25").Value = "12.89"
$ws.Range("E25").Value = "  -3.28%  "
$ws.Range("E26").Value = "  -2.16%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").Value = "8.99"
$ws.Range("E28").Value = "  +8.23%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "2.66"
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("E31").Value = "  -2.27%  "
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("E33").Value = "  -1.62%  "
$ws.Range("D34").Value = "26.61"
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("D35").Value = "0.0₃0849"
$ws.Range("E35").Value = "  -3.06%  "
$ws.Range("E36").Value = "  +2.66%  "
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("D38").Value = "2.29"
$ws.Range("E38").Value = "  -4.17%  "
$ws.Range("D39").Value = "5.98"
$ws.Range("E39").Value = "  -1.35%  "
$ws.Range("D40").Value = "50.32"
$ws.Range("E40").Value = "  -1.30%  "
$ws.Range("D41").Value = "433.64"
$ws.Range("E41").Value = "  -0.49%  "
$ws.Range("D42").Value = "8.69"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").Value = "2.878.07"
$ws.Range("E44").Value = "  -2.06%  "
$ws.Range("D45").Value = "0.269"
$ws.Range("E45").Value = "  -3.48%  "
$ws.Range("E46").Value = "  -3.61%  "
$ws.Range("D47").Value = "35.76"
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").Value = "123.37"
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").Value = "0.109"
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("D51").Value = "24.04"
$ws.Range("E51").Value = "  -2.58%  "
